# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# Adds/normalises a "schema" row (row 3) on the Body / response-code sheets so
# each points at its named schema object (Request/Response/errorResponse
# variants), replacing the old ad-hoc "dateTime" breakdown rows and dropping
# the now-redundant follow-up rows that described the individual fields.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow3($SheetName, $Col_A, $SchemaName) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range("A3").Value = $Col_A
    $ws.Range("B3").Value = $SchemaName
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $SchemaName
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# --- Body: replace the dateTime/settlementBIC/init/end breakdown (rows 3-6)
# with a single pointer to the request schema. Rows below are cleared (not
# deleted) so unrelated absolute-ish ranges (conditional formatting, data
# validation) that spill into the lower part of the sheet keep their
# original extents. ---
$wsBody = $wb.Worksheets.Item("Body")
$wsBody.Range("A4:O6").ClearContents()
Set-SchemaRow3 "Body" "body" "revokeChangeSettlementBIC.230216Request"

# --- 200: replace the dateTime/commandRef breakdown (rows 3-4) with a
# single pointer to the response schema. ---
$ws200 = $wb.Worksheets.Item("200")
$ws200.Range("A4:O4").ClearContents()
Set-SchemaRow3 "200" "content" "revokeChangeSettlementBIC.230216Response"

# --- 204: no existing breakdown rows; just add the pointer row to the same
# response schema as 200. ---
Set-SchemaRow3 "204" "content" "revokeChangeSettlementBIC.230216Response"

# --- 400: replace the dateTime/errorCode/errorCodeDescription/requestId
# breakdown (rows 3-6) with a single pointer to the generic error schema. ---
$ws400 = $wb.Worksheets.Item("400")
$ws400.Range("A4:O6").ClearContents()
Set-SchemaRow3 "400" "content" "errorResponse"

# --- 401, 403, 404, 429, 500: no existing breakdown rows; add the pointer
# row to the shared errorResponse1 schema. ---
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    Set-SchemaRow3 $sheetName "content" "errorResponse1"
}
